$d = $word.ActiveDocument

$replacements = @(
    @("315×4=", "653×7="),
    @("371×8=", "950×4="),
    @("147×2=", "473×9="),
    @("133×6=", "696×7="),
    @("700×8=", "774×5="),
    @("776×7=", "515×2="),
    @("712×4=", "631×5="),
    @("980×8=", "482×7="),
    @("541×6=", "465×3="),
    @("719×3=", "118×7="),
    @("576×3=", "219×9="),
    @("142×8=", "890×4="),
    @("122×5=", "981×8="),
    @("604×9=", "188×8="),
    @("263×2=", "610×7="),
    @("344×3=", "792×6="),
    @("304×4=", "651×2="),
    @("832×4=", "733×6="),
    @("132×7=", "148×4="),
    @("745×9=", "114×7="),
    @("882×3=", "556×4="),
    @("388×4=", "233×9="),
    @("147×3=", "833×5="),
    @("104×8=", "753×7="),
    @("459×4=", "221×7="),
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
